$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '41.058.41'
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = '  -1.71%  '

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '2.170.65'
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = '  -2.60%  '

$ws.Range("E4").Value = '  -0.05%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '247.48'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -1.06%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '0.611'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -2.79%  '

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '65.50'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  -8.38%  '

$ws.Range("E8").Value = '  -0.04%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.560'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  -4.35%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '59.67'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  +2.36%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.0923'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  -4.83%  '

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '35.33'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  -14.61%  '

$ws.Range("E13").Value = '  -2.25%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '6.82'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  -2.69%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '2.490.75'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  -2.72%  '

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '14.17'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  -5.38%  '

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '0.846'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  -2.13%  '

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '2.168.69'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  -2.73%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '41.046.21'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  -1.71%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '0.0₃0933'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  -3.89%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '71.18'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  -2.25%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '6.05'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  -2.64%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '228.50'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  -2.93%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '2.07'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  -4.04%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '3.83'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  -9.34%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '1.00'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  +0.07%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '11.17'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  +6.03%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '2.41'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  -4.90%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '3.72'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  -5.74%  '

$ws.Range("E30").Value = '  -3.47%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '167.49'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  -2.28%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '20.13'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  -2.95%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '0.121'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  -0.80%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '5.62'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  -0.52%  '

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.0739'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  +2.24%  '

$ws.Range("E36").Value = '  -3.65%  '

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '4.52'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  -3.85%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '3.97'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  +1.07%  '

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '24.04'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  -7.98%  '

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.0300'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  -0.97%  '

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '2.18'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  -5.04%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '5.43'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  -8.79%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '4.84'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  -0.78%  '

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '60.30'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  -11.81%  '

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '11.09'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  -6.41%  '

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.190'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  -8.73%  '

$ws.Range("B47").Value = 'FraxShare'
$ws.Range("C47").Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '8.43'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  -4.48%  '

$ws.Range("B48").Value = 'BinanceUSD'
$ws.Range("C48").Value = 'https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd'
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '1.00'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  -0.24%  '

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '0.0986'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  -3.41%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '1.14'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  -1.83%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '1.14'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  -4.72%  '
